# Update "Analise Financeira PROJ1" – Projeto 1 analysis sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CÁLCULO - VPL-TIR-PAYBACK")

# 1. Update the TMA (Taxa Minima de Atratividade) input from 10% to 11%.
#    All dependent formulas (VPL, payback, etc.) recalculate automatically.
$ws.Range("G8").Value = 0.11

# 2. Update the narrative analysis text in the merged cell F14:G19.
$ws.Range("F14").Value = "Com uma Taxa Mínima de Atratividade (TMA) de 11%, o projeto continua sendo financeiramente viável. O Valor Presente Líquido (VPL) é positivo, atingindo R`$ 266.545,94, o que indica que o valor presente dos fluxos de caixa futuros supera o investimento inicial de R`$ 1.584.000,00. O período de payback é de aproximadamente 5,29 anos, destacando um tempo relativamente longo para recuperar o capital investido. A Taxa Interna de Retorno (TIR) de 15,68% supera a TMA de 11%, sugerindo que o projeto oferece um retorno superior ao custo de oportunidade do capital. Com um Retorno sobre o Investimento (ROI) de 77%, o projeto demonstra uma boa relação custo-benefício e reforça sua atratividade financeira. Portanto, é recomendada sua aprovação, pois o projeto oferece um retorno significativo com risco controlado."

# 3. The longer text needs a taller row to display fully (row 19 holds the
#    merged range's bottom edge / custom height).
$ws.Rows.Item(19).RowHeight = 110.4

# 4. Update the sheet's active selection, matching where the author ended up.
[void]$ws.Range("K19").Select()
